# Auto-generated Excel COM-interop script
# Applies value updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets
# as produced by the scheduled Sheets runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value2 = 3879.1177
$ws.Cells.Item(5, 9).Value2 = 1330.5
$ws.Cells.Item(5, 11).Value2 = 1330.5
$ws.Cells.Item(5, 13).Value2 = -1215.5

$ws.Cells.Item(28, 8).Value2 = 1801.7693
$ws.Cells.Item(28, 9).Value2 = 920.6667
$ws.Cells.Item(28, 11).Value2 = 920.6667
$ws.Cells.Item(28, 13).Value2 = -435.6667

$ws.Cells.Item(43, 8).Value2 = 3379.4666
$ws.Cells.Item(43, 9).Value2 = 2300.3333
$ws.Cells.Item(43, 11).Value2 = 2300.3333
$ws.Cells.Item(43, 13).Value2 = -2231.3333

$ws.Cells.Item(80, 8).Value2 = 2088
$ws.Cells.Item(80, 9).Value2 = 2300.6667
$ws.Cells.Item(80, 10).Value2 = 1928.5
$ws.Cells.Item(80, 11).Value2 = 6902.000100000001
$ws.Cells.Item(80, 12).Value2 = 5785.5
$ws.Cells.Item(80, 13).Value2 = -5904.000100000001
$ws.Cells.Item(80, 14).Value2 = -7781.5

$ws.Cells.Item(83, 8).Value2 = 2088
$ws.Cells.Item(83, 9).Value2 = 2300.6667
$ws.Cells.Item(83, 10).Value2 = 1928.5
$ws.Cells.Item(83, 11).Value2 = 20706.0003
$ws.Cells.Item(83, 12).Value2 = 17356.5
$ws.Cells.Item(83, 13).Value2 = -15714.0003
$ws.Cells.Item(83, 14).Value2 = -27340.5

$ws.Cells.Item(138, 8).Value2 = 4781.517
$ws.Cells.Item(138, 9).Value2 = 1396
$ws.Cells.Item(138, 10).Value2 = 5323.2
$ws.Cells.Item(138, 11).Value2 = 4188
$ws.Cells.Item(138, 12).Value2 = 15969.6
$ws.Cells.Item(138, 13).Value2 = 952
$ws.Cells.Item(138, 14).Value2 = -26249.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value2 = 495.875
$ws.Cells.Item(4, 9).Value2 = 473.16666
$ws.Cells.Item(4, 11).Value2 = 473.16666
$ws.Cells.Item(4, 13).Value2 = -357.16666

$ws.Cells.Item(32, 8).Value2 = 1917067.2
$ws.Cells.Item(32, 9).Value2 = 619.30615
$ws.Cells.Item(32, 10).Value2 = 12351061
$ws.Cells.Item(32, 11).Value2 = 619.30615
$ws.Cells.Item(32, 12).Value2 = 12351061
$ws.Cells.Item(32, 13).Value2 = -332.30615
$ws.Cells.Item(32, 14).Value2 = -12351635

$ws.Cells.Item(88, 8).Value2 = 5492.2
$ws.Cells.Item(88, 9).Value2 = 1768.25
$ws.Cells.Item(88, 10).Value2 = 9748.143
$ws.Cells.Item(88, 11).Value2 = 1768.25
$ws.Cells.Item(88, 12).Value2 = 9748.143
$ws.Cells.Item(88, 13).Value2 = -1362.25
$ws.Cells.Item(88, 14).Value2 = -10560.143

$ws.Cells.Item(91, 8).Value2 = 5492.2
$ws.Cells.Item(91, 9).Value2 = 1768.25
$ws.Cells.Item(91, 10).Value2 = 9748.143
$ws.Cells.Item(91, 11).Value2 = 1768.25
$ws.Cells.Item(91, 12).Value2 = 9748.143
$ws.Cells.Item(91, 13).Value2 = -364.25
$ws.Cells.Item(91, 14).Value2 = -12556.143

$ws.Cells.Item(132, 8).Value2 = 1079766.1
$ws.Cells.Item(132, 9).Value2 = 1445534.6
$ws.Cells.Item(132, 10).Value2 = 119623.875
$ws.Cells.Item(132, 11).Value2 = 4336603.800000001
$ws.Cells.Item(132, 12).Value2 = 358871.625
$ws.Cells.Item(132, 13).Value2 = -4334073.800000001
$ws.Cells.Item(132, 14).Value2 = -363931.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value2 = 5120
$ws.Cells.Item(64, 10).Value2 = 9883
$ws.Cells.Item(64, 12).Value2 = 9883
$ws.Cells.Item(64, 14).Value2 = -10333

$ws.Cells.Item(67, 8).Value2 = 5120
$ws.Cells.Item(67, 10).Value2 = 9883
$ws.Cells.Item(67, 12).Value2 = 9883
$ws.Cells.Item(67, 14).Value2 = -11443

$ws.Cells.Item(104, 8).Value2 = 0
$ws.Cells.Item(104, 10).Value2 = 0
$ws.Cells.Item(104, 12).Value2 = 0
$ws.Cells.Item(104, 14).ClearContents()

$ws.Cells.Item(106, 8).Value2 = 21417.5
$ws.Cells.Item(106, 10).Value2 = 21417.5
$ws.Cells.Item(106, 12).Value2 = 21417.5
$ws.Cells.Item(106, 14).Value2 = -23941.5

$ws.Cells.Item(110, 8).Value2 = 34899
$ws.Cells.Item(110, 10).Value2 = 34899
$ws.Cells.Item(110, 12).Value2 = 34899
$ws.Cells.Item(110, 14).Value2 = -43079

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value2 = 346.07693
$ws.Cells.Item(7, 9).Value2 = 290.1
$ws.Cells.Item(7, 11).Value2 = 290.1
$ws.Cells.Item(7, 13).Value2 = -177.1

$ws.Cells.Item(62, 8).Value2 = 17375
$ws.Cells.Item(62, 10).Value2 = 17375
$ws.Cells.Item(62, 12).Value2 = 17375
$ws.Cells.Item(62, 14).Value2 = -18623

$ws.Cells.Item(65, 8).Value2 = 17375
$ws.Cells.Item(65, 10).Value2 = 17375
$ws.Cells.Item(65, 12).Value2 = 86875
$ws.Cells.Item(65, 14).Value2 = -93115

$ws.Cells.Item(86, 8).Value2 = 11808.056
$ws.Cells.Item(86, 9).Value2 = 13269.667
$ws.Cells.Item(86, 10).Value2 = 11077.25
$ws.Cells.Item(86, 11).Value2 = 13269.667
$ws.Cells.Item(86, 12).Value2 = 11077.25
$ws.Cells.Item(86, 13).Value2 = -12146.667
$ws.Cells.Item(86, 14).Value2 = -13323.25

$ws.Cells.Item(89, 8).Value2 = 11808.056
$ws.Cells.Item(89, 9).Value2 = 13269.667
$ws.Cells.Item(89, 10).Value2 = 11077.25
$ws.Cells.Item(89, 11).Value2 = 66348.33499999999
$ws.Cells.Item(89, 12).Value2 = 55386.25
$ws.Cells.Item(89, 13).Value2 = -60732.33499999999
$ws.Cells.Item(89, 14).Value2 = -66618.25

$ws.Cells.Item(94, 8).Value2 = 62504252
$ws.Cells.Item(94, 9).Value2 = 166668510
$ws.Cells.Item(94, 10).Value2 = 5698.1
$ws.Cells.Item(94, 11).Value2 = 166668510
$ws.Cells.Item(94, 12).Value2 = 5698.1
$ws.Cells.Item(94, 13).Value2 = -166668059
$ws.Cells.Item(94, 14).Value2 = -6600.1

$ws.Cells.Item(132, 8).Value2 = 3834.8108
$ws.Cells.Item(132, 9).Value2 = 2944.8215
$ws.Cells.Item(132, 11).Value2 = 8834.4645
$ws.Cells.Item(132, 13).Value2 = -6304.4645

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value2 = 3454.7144
$ws.Cells.Item(11, 9).Value2 = 3437.8
$ws.Cells.Item(11, 11).Value2 = 10313.4
$ws.Cells.Item(11, 13).Value2 = -10173.4

$ws.Cells.Item(37, 8).Value2 = 218120.83
$ws.Cells.Item(37, 10).Value2 = 218120.83
$ws.Cells.Item(37, 12).Value2 = 654362.49
$ws.Cells.Item(37, 14).Value2 = -654586.49

$ws.Cells.Item(80, 8).Value2 = 4200
$ws.Cells.Item(80, 9).Value2 = 2000
$ws.Cells.Item(80, 10).Value2 = 4750
$ws.Cells.Item(80, 11).Value2 = 6000
$ws.Cells.Item(80, 12).Value2 = 14250
$ws.Cells.Item(80, 13).Value2 = -5064
$ws.Cells.Item(80, 14).Value2 = -16122

$ws.Cells.Item(83, 8).Value2 = 4200
$ws.Cells.Item(83, 9).Value2 = 2000
$ws.Cells.Item(83, 10).Value2 = 4750
$ws.Cells.Item(83, 11).Value2 = 18000
$ws.Cells.Item(83, 12).Value2 = 42750
$ws.Cells.Item(83, 13).Value2 = -13320
$ws.Cells.Item(83, 14).Value2 = -52110

$ws.Cells.Item(87, 8).Value2 = 20688.715
$ws.Cells.Item(87, 9).Value2 = 5273.6665
$ws.Cells.Item(87, 11).Value2 = 15820.9995
$ws.Cells.Item(87, 13).Value2 = -14572.9995

$ws.Cells.Item(90, 8).Value2 = 20688.715
$ws.Cells.Item(90, 9).Value2 = 5273.6665
$ws.Cells.Item(90, 11).Value2 = 47462.9985
$ws.Cells.Item(90, 13).Value2 = -41222.9985

$ws.Cells.Item(134, 8).Value2 = 4816.7744
$ws.Cells.Item(134, 10).Value2 = 0
$ws.Cells.Item(134, 12).Value2 = 0
$ws.Cells.Item(134, 14).ClearContents()

$ws.Cells.Item(140, 8).Value2 = 75762310
$ws.Cells.Item(140, 9).Value2 = 119052220
$ws.Cells.Item(140, 10).Value2 = 4973.75
$ws.Cells.Item(140, 11).Value2 = 357156660
$ws.Cells.Item(140, 12).Value2 = 14921.25
$ws.Cells.Item(140, 13).Value2 = -357151480
$ws.Cells.Item(140, 14).Value2 = -25281.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value2 = 9830.200000000001
$ws.Cells.Item(113, 9).Value2 = 2075.75
$ws.Cells.Item(113, 11).Value2 = 2075.75
$ws.Cells.Item(113, 13).Value2 = 94.25

$ws.Cells.Item(132, 8).Value2 = 66671560
$ws.Cells.Item(132, 9).Value2 = 125004260
$ws.Cells.Item(132, 10).Value2 = 5620.7144
$ws.Cells.Item(132, 11).Value2 = 375012780
$ws.Cells.Item(132, 12).Value2 = 16862.1432
$ws.Cells.Item(132, 13).Value2 = -375010250
$ws.Cells.Item(132, 14).Value2 = -21922.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value2 = 21740978
$ws.Cells.Item(46, 9).Value2 = 800.1
$ws.Cells.Item(46, 10).Value2 = 38464190
$ws.Cells.Item(46, 11).Value2 = 800.1
$ws.Cells.Item(46, 12).Value2 = 38464190
$ws.Cells.Item(46, 13).Value2 = -612.1
$ws.Cells.Item(46, 14).Value2 = -38464566

$ws.Cells.Item(93, 8).Value2 = 1271.4828
$ws.Cells.Item(93, 9).Value2 = 1437
$ws.Cells.Item(93, 11).Value2 = 1437
$ws.Cells.Item(93, 13).Value2 = -189

$ws.Cells.Item(132, 8).Value2 = 2763.3
$ws.Cells.Item(132, 9).Value2 = 1830.6666
$ws.Cells.Item(132, 11).Value2 = 5491.9998
$ws.Cells.Item(132, 13).Value2 = -2961.9998

$ws.Cells.Item(136, 8).Value2 = 23812180
$ws.Cells.Item(136, 9).Value2 = 62501050
$ws.Cells.Item(136, 11).Value2 = 187503150
$ws.Cells.Item(136, 13).Value2 = -187500600
